$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Voltooide coachings")
$tbl = $ws.ListObjects.Item("Tabel3")

# Insert a new column into the table just before the last column ("Ingevoerd in SSF")
$lastCol = $tbl.ListColumns.Item($tbl.ListColumns.Count)
$newCol = $tbl.ListColumns.Add($lastCol.Index)
$newCol.Name = "Gesprek"

# Update a data value: row with P-nr 33639 (row 42) "Hercoaching noodzakelijk" -> "ja"
$ws.Cells.Item(42, 16).Value = "ja"
